# DATA: update 26 Maret 2020
# Updates the last data row (A38 = 2020-03-25) with the refreshed figures
# from the 26 Maret 2020 site pull, and drops the now-stale "catatan"
# footnote that used to sit on that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 38 (tanggal 2020-03-25) figure refresh -----------------------
$ws.Cells.Item(38, 2).Value = 3822   # B38 jumlah_periksa   3332 -> 3822
$ws.Cells.Item(38, 6).Value = 3032   # F38 negatif          2625 -> 3032
$ws.Cells.Item(38, 7).Value = 34     # G38 proses_periksa      0 -> 34
$ws.Cells.Item(38, 8).Value = 0      # H38 kasus_perawatan   701 -> 0

# The "(informasi diambil di situs, pada 2020-03-25 21:06 WIB)" note in
# I38 no longer applies now that the row has been refreshed - clear it
# (this also drops the now-unused shared string).
$ws.Cells.Item(38, 9).ClearContents()

# --- Selection / view bookkeeping (matches the saved view state) ------
$ws.Range("F37").Select()
